$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 482, shifting existing rows 482-524 down to 484-526
$ws.Range("A482:A483").EntireRow.Insert()

# Populate new row 482 with the new record (date 45166, Volumen 300, ...)
$ws.Range("A482").Value = 8
$ws.Range('B482').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C482').Value = 'Coquimbo'
$ws.Range("D482").Value = 45166
$ws.Range("E482").Value = 4
$ws.Range("F482").Value = 100112003
$ws.Range('G482').Value = 'Ajo'
$ws.Range('H482').Value = 'Chino'
$ws.Range('I482').Value = 'Primera'
$ws.Range("J482").Value = 300
$ws.Range("K482").Value = 21000
$ws.Range("L482").Value = 22000
$ws.Range("M482").Value = 21500
$ws.Range('N482').Value = '$/caja 10 kilos'
$ws.Range('O482').Value = 'China'
$ws.Range("P482").Value = 2150
$ws.Range("Q482").Value = 10
$ws.Range('R482').Value = 'Hortaliza'

# Populate new row 483 with the new record (date 45166, Volumen 340, ...)
$ws.Range("A483").Value = 8
$ws.Range('B483').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C483').Value = 'Coquimbo'
$ws.Range("D483").Value = 45166
$ws.Range("E483").Value = 4
$ws.Range("F483").Value = 100112003
$ws.Range('G483').Value = 'Ajo'
$ws.Range('H483').Value = 'Chino'
$ws.Range('I483').Value = 'Primera'
$ws.Range("J483").Value = 340
$ws.Range("K483").Value = 23000
$ws.Range("L483").Value = 24000
$ws.Range("M483").Value = 23500
$ws.Range('N483').Value = '$/malla 10 kilos'
$ws.Range('O483').Value = 'China'
$ws.Range("P483").Value = 2350
$ws.Range("Q483").Value = 10
$ws.Range('R483').Value = 'Hortaliza'
